$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.498.83'
$ws.Range("E2").Value = '  -0.94%  '

$ws.Range("D3").Value = '1.694.73'
$ws.Range("E3").Value = '  -0.39%  '

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.41%  '

$ws.Range("D5").Value = "'315.78"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.22%  '

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -0.39%  '

$ws.Range("D7").Value = "'0.3919"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -0.62%  '

$ws.Range("D8").Value = "'0.4063"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +0.53%  '

$ws.Range("E9").Value = '  -2.79%  '

$ws.Range("E10").Value = '  -0.35%  '

$ws.Range("D11").Value = "'52.57"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -2.07%  '

$ws.Range("D12").Value = "'0.08788"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -0.97%  '

$ws.Range("D13").Value = "'26.86"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +13.38%  '

$ws.Range("D14").Value = "'7.514"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("D15").Value = "'8.132"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("D16").Value = "'0.00001349"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +1.72%  '

$ws.Range("D17").Value = '1.691.04'
$ws.Range("E17").Value = '  -0.97%  '

$ws.Range("E18").Value = '  -1.74%  '

$ws.Range("D19").Value = "'0.07165"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +1.65%  '

$ws.Range("D20").Value = "'20.61"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +4.47%  '

$ws.Range("D21").Value = "'7.289"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +3.03%  '

$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("E23").Value = '  -2.17%  '

$ws.Range("D24").Value = '24.490.94'
$ws.Range("E24").Value = '  -0.95%  '

$ws.Range("D25").Value = "'3.033"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -4.56%  '

$ws.Range("D26").Value = "'2.325"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -1.93%  '

$ws.Range("D27").Value = "'22.72"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -0.23%  '

$ws.Range("D28").Value = "'167.17"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +3.01%  '

$ws.Range("D29").Value = "'8.556"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -1.68%  '

$ws.Range("D30").Value = "'5.401"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +4.28%  '

$ws.Range("D31").Value = "'140.09"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +3.10%  '

$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").Value = "'2.216"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +11.24%  '

$ws.Range("B33").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C33").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D33").Value = '1.876.09'
$ws.Range("E33").Value = '  -0.97%  '

$ws.Range("D34").Value = "'0.08771"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -3.22%  '

$ws.Range("D35").Value = "'7.298"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -4.78%  '

$ws.Range("E36").Value = '  -3.57%  '

$ws.Range("D37").Value = "'0.03043"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +9.48%  '

$ws.Range("D38").Value = "'0.2801"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +1.55%  '

$ws.Range("D39").Value = "'10.93"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -1.74%  '

$ws.Range("D40").Value = "'0.09177"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +0.07%  '

$ws.Range("E41").Value = '  -2.37%  '

$ws.Range("D42").Value = "'0.8029"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +4.20%  '

$ws.Range("D43").Value = "'1.484"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +1.47%  '

$ws.Range("D44").Value = "'17.60"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +9.68%  '

$ws.Range("D45").Value = "'2.668"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +3.57%  '

$ws.Range("D46").Value = "'0.7281"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +1.17%  '

$ws.Range("D47").Value = "'4.268"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +1.13%  '

$ws.Range("D48").Value = "'1.408"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +4.55%  '

$ws.Range("D49").Value = "'1.000"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -0.38%  '

$ws.Range("D50").Value = "'141.16"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +0.39%  '

$ws.Range("D51").Value = "'0.08158"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +2.16%  '
